$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellValue($addr, $val) {
    $ws.Range($addr).Value = $val
}

# --- Row 2 ---
Set-CellValue "A2" 'https://openalex.org/W4324031766'
Set-CellValue "B2" 'Genuine Explanation and the Strong Minimalist Thesis'
Set-CellValue "C2" 'Genuine Explanation and the Strong Minimalist Thesis'
Set-CellValue "E2" 'Abstract The goal of theoretical inquiry is explanation: Why this, and not that? In the study of language, search for explanatory theory proceeds at two levels: for individual languages (a generative grammar in the broad sense) and for the general faculty of language fl ( ug ), the latter apparently a true species property, common to humans and without significant analogue in the animal world. ug must meet several conditions: learnability, evolvability, coverage. These conditions appear to conflict, and are far more severe than had earlier been supposed. A solution to the conundrum would be satisfaction of smt for ug combined with recourse to language-independent principles of computational efficiency, with diversity sequestered in components of language subject to simple algorithmic search. For the first time, hopes for such an outcome seem to be on the horizon, with significant implications if the hopes can be realized. I will outline some current work on these topics.'
Set-CellValue "F2" '2023-01-25'
Set-CellValue "G2" 'Cognitive Semantics'
Set-CellValue "H2" 'https://openalex.org/S4210197614'
Set-CellValue "J2" '2352-6408'
Set-CellValue "K2" 'https://doi.org/10.1163/23526416-bja10040'
Set-CellValue "L2" 'https://brill.com/downloadpdf/journals/cose/8/3/article-p347_002.pdf'
Set-CellValue "N2" 'publishedVersion'
Set-CellValue "O2" '347'
Set-CellValue "P2" '365'
Set-CellValue "Q2" '8'
Set-CellValue "S2" $true
Set-CellValue "T2" $true
Set-CellValue "U2" 'bronze'
Set-CellValue "V2" 'https://brill.com/downloadpdf/journals/cose/8/3/article-p347_002.pdf'
Set-CellValue "Z2" 8
Set-CellValue "AB2" 2023
Set-CellValue "AC2" 'https://api.openalex.org/works?filter=cites:W4324031766'
Set-CellValue "AE2" 'https://doi.org/10.1163/23526416-bja10040'

# --- Row 3 ---
Set-CellValue "A3" 'https://openalex.org/W4388460035'
Set-CellValue "B3" 'Consciousness Is Quantum State Reduction Which Creates the Flow of Time'
Set-CellValue "C3" 'Consciousness Is Quantum State Reduction Which Creates the Flow of Time'
Set-CellValue "E3" 'Abstract In neuroscience, the flow of time is a conscious experience produced by the brain. But in physics, time is either a process, or a dimension in four-dimensional spacetime geometry. Could all three explanations be correct? The Penrose–Hameroff ‘Orch OR’ theory suggests consciousness is a sequence of discrete, irreversible quantum state reductions occurring at an objective threshold (‘objective reduction’, ‘OR’) a process in fundamental spacetime geometry. These ‘self-collapses’ of the quantum wavefunction are ‘orchestrated’ in microtubules inside brain neurons. Each Orch OR event selects microtubule states which purposefully regulate neuronal functions, and provide sequences of ‘NOW’ moments of conscious experience. Connected to fundamental spacetime geometry by Penrose ‘OR’, consciousness is quantum state reduction, a set of irreversible steps which ‘ratchet forward’ in the fine scale geometry of the universe, creating a flow of time.'
Set-CellValue "F3" '2023-11-06'
Set-CellValue "G3" 'Timing & Time Perception'
Set-CellValue "H3" 'https://openalex.org/S4210181418'
Set-CellValue "J3" '2213-445X'
Set-CellValue "K3" 'https://doi.org/10.1163/22134468-bja10098'
Set-CellValue "L3" 'https://brill.com/downloadpdf/view/journals/time/aop/article-10.1163-22134468-bja10098/article-10.1163-22134468-bja10098.pdf'
Set-CellValue "M3" 'cc-by'
Set-CellValue "N3" 'publishedVersion'
Set-CellValue "O3" '158'
Set-CellValue "P3" '167'
Set-CellValue "Q3" '12'
Set-CellValue "S3" $true
Set-CellValue "T3" $true
Set-CellValue "U3" 'hybrid'
Set-CellValue "V3" 'https://brill.com/downloadpdf/view/journals/time/aop/article-10.1163-22134468-bja10098/article-10.1163-22134468-bja10098.pdf'
Set-CellValue "Z3" 3
Set-CellValue "AB3" 2023
Set-CellValue "AC3" 'https://api.openalex.org/works?filter=cites:W4388460035'
Set-CellValue "AE3" 'https://doi.org/10.1163/22134468-bja10098'

# --- Row 4 ---
Set-CellValue "A4" 'https://openalex.org/W4313596414'
Set-CellValue "B4" 'Mummy Labels: A Witness to the Use and Processing of Wood in Roman Egypt'
Set-CellValue "C4" 'Mummy Labels: A Witness to the Use and Processing of Wood in Roman Egypt'
Set-CellValue "E4" 'Abstract Mummy labels are relics found in large quantities in Egypt, often in an excellent state of preservation (like most woods preserved in arid environments). As a result, they are widespread in Roman Egyptian collections of many museums. These labels reflect funerary practices that possess Egyptian and Roman influences and are an important source of historical and archaeological information. These corpora of mummy labels offer several possibilities for investigation. The inscriptions on these labels have been the subject of an international project (Death on the Nile) in which all accessible objects were recorded in a database. However, the potential of these funerary objects extend beyond the inscriptions to the methods of manufacturing and cutting, the choice of species used, and their dendrochronological potential to better define their chronology and possibly their provenance. The study of mummy labels allows us to propose a new typology, some forms of which seem to be limited to certain necropolises. Mummy labels, whether made by the family of the deceased or by specific workshops, show that their realizations vary greatly, ranging from coarse specimens to others with beautiful detailing. They are made from endemic as well as imported species, which are symbolic of long-distance trade, especially for conifer trees, which are well represented. Their dendrochronological potential has also been demonstrated in numerous studies, some of which have allowed the identification of labels from the same tree, supported by inscriptions attesting to the same family relationship.'
Set-CellValue "F4" '2023-01-04'
Set-CellValue "G4" 'International Journal of Wood Culture'
Set-CellValue "H4" 'https://openalex.org/S4210185514'
Set-CellValue "J4" '2772-3186'
Set-CellValue "K4" 'https://doi.org/10.1163/27723194-bja10017'
Set-CellValue "L4" 'https://brill.com/downloadpdf/journals/ijwc/aop/article-10.1163-27723194-bja10017/article-10.1163-27723194-bja10017.pdf'
Set-CellValue "M4" 'cc-by'
Set-CellValue "O4" '192'
Set-CellValue "P4" '223'
Set-CellValue "Q4" '3'
Set-CellValue "R4" '1-3'
Set-CellValue "U4" 'diamond'
Set-CellValue "V4" 'https://brill.com/downloadpdf/journals/ijwc/aop/article-10.1163-27723194-bja10017/article-10.1163-27723194-bja10017.pdf'
Set-CellValue "Z4" 3
Set-CellValue "AB4" 2023
Set-CellValue "AC4" 'https://api.openalex.org/works?filter=cites:W4313596414'
Set-CellValue "AE4" 'https://doi.org/10.1163/27723194-bja10017'

# --- Row 5 ---
Set-CellValue "A5" 'https://openalex.org/W4385310714'
Set-CellValue "B5" 'The Virtues of Sustainability, edited by Jason Kawall'
Set-CellValue "C5" 'The Virtues of Sustainability, edited by Jason Kawall'
Set-CellValue "F5" '2023-07-24'
Set-CellValue "G5" 'Journal of Moral Philosophy'
Set-CellValue "H5" 'https://openalex.org/S96509893'
Set-CellValue "J5" '1740-4681'
Set-CellValue "K5" 'https://doi.org/10.1163/17455243-20030008'
Set-CellValue "O5" '362'
Set-CellValue "P5" '365'
Set-CellValue "Q5" '20'
Set-CellValue "R5" '3-4'
Set-CellValue "AB5" 2023
Set-CellValue "AC5" 'https://api.openalex.org/works?filter=cites:W4385310714'
Set-CellValue "AE5" 'https://doi.org/10.1163/17455243-20030008'

# --- Row 6 ---
Set-CellValue "A6" 'https://openalex.org/W4386686640'
Set-CellValue "B6" 'Homoerotic and Homosexual Perspectives in Medieval Poetry and Verse Narratives: Indirect Evidence of a Hidden Discourse'
Set-CellValue "C6" 'Homoerotic and Homosexual Perspectives in Medieval Poetry and Verse Narratives: Indirect Evidence of a Hidden Discourse'
Set-CellValue "E6" 'Abstract Although it proves to be a difficult task, we still can identify more literary texts from the Middle Ages addressing homoerotic love than we might have expected. Even when poets voiced severe criticism and radically condemned homosexuality, their comments serve us well to identify more specifically the actual discourse behind the official scene. Although legal and Church authorities consistently characterized ‘sodomy’ as one of the worst sins a Christian could commit, since late antiquity, and certainly throughout the Middle Ages, the phenomenon, a biological fact, existed, of course, and was also addressed in veiled or open language. This article examines a selection of relevant literary and didactic works that shed more light on this issue.'
Set-CellValue "F6" '2023-09-06'
Set-CellValue "G6" 'Amsterdamer Beiträge zur älteren Germanistik'
Set-CellValue "H6" 'https://openalex.org/S4210206867'
Set-CellValue "J6" '0165-7305'
Set-CellValue "K6" 'https://doi.org/10.1163/18756719-12340294'
Set-CellValue "O6" '234'
Set-CellValue "P6" '249'
Set-CellValue "Q6" '83'
Set-CellValue "AB6" 2023
Set-CellValue "AC6" 'https://api.openalex.org/works?filter=cites:W4386686640'
Set-CellValue "AE6" 'https://doi.org/10.1163/18756719-12340294'

# --- Row 7 ---
Set-CellValue "A7" 'https://openalex.org/W4388943180'
Set-CellValue "B7" 'Doctrinal Engagements and Disengagements: Yongming Yanshou and His Legacies'
Set-CellValue "C7" 'Doctrinal Engagements and Disengagements: Yongming Yanshou and His Legacies'
Set-CellValue "E7" 'Abstract This paper looks at the vexed relationship of doctrine, or teaching (C. jiao /K. kyo /J. kyō 教 ) in the three kindred traditions subsumed under the rubric of the Sino-East Asian graph 禪 , known through their distinctive pronunciations in modern languages as Chan, Sŏn, and Zen. While the stipulation of these traditions as ‘a special/separate transmission outside the teachings’; (jiaowai biechuan 教外別傳 ) presumes independence from Mahayana doctrinal teachings, the reality, as we know, was much more complicated. In this paper, I use Yongming Yanshou 永明延壽 (904–975/6), one of the most prominent Chan figures to promote doctrinal engagement, as a barometer to look at how doctrinal engagements and disengagements are regarded throughout each tradition. Perspectives on Yanshou, a figure at once revered and marginalized, unlock key features of each of these three interconnected traditions, what they share and how they disagree. Fundamentally, perspectives on doctrinal engagements and disengagements are rooted in seminal Chan disputes over the nature and value of Buddhist teaching, and Yanshou is a conduit for these disputes. Given the theme of the conference, ‘How Zen Became Chan’; I also look at the discrepancies these disputes reveal between modern Rinzai Zen orthodoxy’s defining of Zen in the modern world and the practice of Chan in China and Sŏn in Korea. The options that these discrepancies reveal are indicative of the relevance of doctrinal entanglements and disentanglements to the contemporary Chan, Sŏn, and Zen worlds.'
Set-CellValue "F7" '2023-11-22'
Set-CellValue "G7" 'Journal of Chan Buddhism'
Set-CellValue "H7" 'https://openalex.org/S4210237834'
Set-CellValue "I7" 'Brill'
Set-CellValue "J7" '2589-7160'
Set-CellValue "K7" 'https://doi.org/10.1163/25897179-12340020'
Set-CellValue "O7" '35'
Set-CellValue "P7" '67'
Set-CellValue "Q7" '3'
Set-CellValue "R7" '1-2'
Set-CellValue "S7" $false
Set-CellValue "T7" $false
Set-CellValue "U7" 'closed'
Set-CellValue "W7" $false
Set-CellValue "X7" 'en'
Set-CellValue "Z7" 0
Set-CellValue "AB7" 2023
Set-CellValue "AC7" 'https://api.openalex.org/works?filter=cites:W4388943180'
Set-CellValue "AE7" 'https://doi.org/10.1163/25897179-12340020'
Set-CellValue "AF7" 'article'
Set-CellValue "AI7" $false
Set-CellValue "AJ7" $false

# --- Row 8 ---
Set-CellValue "A8" 'https://openalex.org/W4389207160'
Set-CellValue "B8" 'Art, Allegory, and the Rise of Shiism in Iran, 1487–1565, written by Chad Kia'
Set-CellValue "C8" 'Art, Allegory, and the Rise of Shiism in Iran, 1487–1565, written by Chad Kia'
Set-CellValue "F8" '2023-11-14'
Set-CellValue "G8" 'Shii Studies Review'
Set-CellValue "H8" 'https://openalex.org/S4210224839'
Set-CellValue "I8" 'Brill'
Set-CellValue "J8" '2468-2462'
Set-CellValue "K8" 'https://doi.org/10.1163/24682470-12340096'
Set-CellValue "O8" '409'
Set-CellValue "P8" '412'
Set-CellValue "Q8" '7'
Set-CellValue "R8" '1-2'
Set-CellValue "S8" $false
Set-CellValue "T8" $false
Set-CellValue "U8" 'closed'
Set-CellValue "W8" $false
Set-CellValue "X8" 'en'
Set-CellValue "Z8" 0
Set-CellValue "AB8" 2023
Set-CellValue "AC8" 'https://api.openalex.org/works?filter=cites:W4389207160'
Set-CellValue "AE8" 'https://doi.org/10.1163/24682470-12340096'
Set-CellValue "AF8" 'article'
Set-CellValue "AI8" $false
Set-CellValue "AJ8" $false

